$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 5853
$ws.Range("E2").Value = 5583
$ws.Range("F2").Value = 5466
$ws.Range("G2").Value = 5452
$ws.Range("H2").Value = 5489
$ws.Range("I2").Value = 5668
$ws.Range("J2").Value = 6196
$ws.Range("K2").Value = 6735
$ws.Range("L2").Value = 7262
$ws.Range("M2").Value = 7431
$ws.Range("N2").Value = 7364
$ws.Range("O2").Value = 7246
$ws.Range("P2").Value = 7135
$ws.Range("Q2").Value = 7068
$ws.Range("R2").Value = 6979
$ws.Range("S2").Value = 6993
$ws.Range("T2").Value = 7041
$ws.Range("U2").Value = 7189
$ws.Range("V2").Value = 7501
$ws.Range("W2").Value = 7454
$ws.Range("X2").Value = 7216
$ws.Range("Y2").Value = 6936
$ws.Range("Z2").Value = 6690
$ws.Range("AA2").Value = 6319

$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 5911
$ws.Range("E3").Value = 5643
$ws.Range("F3").Value = 5527
$ws.Range("G3").Value = 5512
$ws.Range("H3").Value = 5545
$ws.Range("I3").Value = 5701
$ws.Range("J3").Value = 6160
$ws.Range("K3").Value = 6630
$ws.Range("L3").Value = 7089
$ws.Range("M3").Value = 7236
$ws.Range("N3").Value = 7148
$ws.Range("O3").Value = 6993
$ws.Range("P3").Value = 6846
$ws.Range("Q3").Value = 6757
$ws.Range("R3").Value = 6640
$ws.Range("S3").Value = 6652
$ws.Range("T3").Value = 6695
$ws.Range("U3").Value = 6827
$ws.Range("V3").Value = 7103
$ws.Range("W3").Value = 7060
$ws.Range("X3").Value = 6842
$ws.Range("Y3").Value = 6585
$ws.Range("Z3").Value = 6359
$ws.Range("AA3").Value = 6020

$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5615
$ws.Range("E4").Value = 5348
$ws.Range("F4").Value = 5232
$ws.Range("G4").Value = 5218
$ws.Range("H4").Value = 5250
$ws.Range("I4").Value = 5406
$ws.Range("J4").Value = 5863
$ws.Range("K4").Value = 6331
$ws.Range("L4").Value = 6788
$ws.Range("M4").Value = 6935
$ws.Range("N4").Value = 6860
$ws.Range("O4").Value = 6726
$ws.Range("P4").Value = 6601
$ws.Range("Q4").Value = 6525
$ws.Range("R4").Value = 6424
$ws.Range("S4").Value = 6434
$ws.Range("T4").Value = 6470
$ws.Range("U4").Value = 6580
$ws.Range("V4").Value = 6812
$ws.Range("W4").Value = 6767
$ws.Range("X4").Value = 6546
$ws.Range("Y4").Value = 6285
$ws.Range("Z4").Value = 6056
$ws.Range("AA4").Value = 5711

$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 5326
$ws.Range("E5").Value = 5060
$ws.Range("F5").Value = 4946
$ws.Range("G5").Value = 4931
$ws.Range("H5").Value = 4964
$ws.Range("I5").Value = 5118
$ws.Range("J5").Value = 5574
$ws.Range("K5").Value = 6040
$ws.Range("L5").Value = 6495
$ws.Range("M5").Value = 6641
$ws.Range("N5").Value = 6578
$ws.Range("O5").Value = 6466
$ws.Range("P5").Value = 6362
$ws.Range("Q5").Value = 6298
$ws.Range("R5").Value = 6214
$ws.Range("S5").Value = 6220
$ws.Range("T5").Value = 6242
$ws.Range("U5").Value = 6308
$ws.Range("V5").Value = 6447
$ws.Range("W5").Value = 6403
$ws.Range("X5").Value = 6186
$ws.Range("Y5").Value = 5928
$ws.Range("Z5").Value = 5702
$ws.Range("AA5").Value = 5362

$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 5038
$ws.Range("E6").Value = 4732
$ws.Range("F6").Value = 4604
$ws.Range("G6").Value = 4590
$ws.Range("H6").Value = 4606
$ws.Range("I6").Value = 4692
$ws.Range("J6").Value = 4804
$ws.Range("K6").Value = 4978
$ws.Range("L6").Value = 5428
$ws.Range("M6").Value = 5859
$ws.Range("N6").Value = 6045
$ws.Range("O6").Value = 6067
$ws.Range("P6").Value = 5974
$ws.Range("Q6").Value = 5915
$ws.Range("R6").Value = 5810
$ws.Range("S6").Value = 5796
$ws.Range("T6").Value = 5771
$ws.Range("U6").Value = 5875
$ws.Range("V6").Value = 5970
$ws.Range("W6").Value = 5804
$ws.Range("X6").Value = 5605
$ws.Range("Y6").Value = 5342
$ws.Range("Z6").Value = 5198
$ws.Range("AA6").Value = 4972

$ws.Range("C7").Value = 8
$ws.Range("D7").Value = 4922
$ws.Range("E7").Value = 4700
$ws.Range("F7").Value = 4627
$ws.Range("G7").Value = 4559
$ws.Range("H7").Value = 4482
$ws.Range("I7").Value = 4555
$ws.Range("J7").Value = 4660
$ws.Range("K7").Value = 4950
$ws.Range("L7").Value = 5358
$ws.Range("M7").Value = 5802
$ws.Range("N7").Value = 6001
$ws.Range("O7").Value = 5977
$ws.Range("P7").Value = 5932
$ws.Range("Q7").Value = 5850
$ws.Range("R7").Value = 5785
$ws.Range("S7").Value = 5790
$ws.Range("T7").Value = 5870
$ws.Range("U7").Value = 6028
$ws.Range("V7").Value = 6129
$ws.Range("W7").Value = 6054
$ws.Range("X7").Value = 5855
$ws.Range("Y7").Value = 5556
$ws.Range("Z7").Value = 5345
$ws.Range("AA7").Value = 5003
